$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up C20:C41 (previously blank/styled placeholder cells) so they
#     mirror the item's Name value in column B, same as every other row.
#     Clear the leftover "text" number-format styling first so the cells
#     come out plain, like the rest of the table. ---
$ws.Range("C20:C41").ClearFormats()
for ($r = 20; $r -le 41; $r++) {
    $nameVal = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $nameVal
}

# --- Correct the swapped Arm slot labels (RightArm1/RightArm2 and
#     LeftArm1/LeftArm2 had been assigned to the wrong armour pieces). ---
$ws.Range("E25").Value = "RightArm2"
$ws.Range("E26").Value = "RightArm1"
$ws.Range("E32").Value = "LeftArm2"
$ws.Range("E33").Value = "LeftArm1"

# --- Add the new Leg item options (rows 42-47). ---
$ws.Range("B42").Value = "Right Leg 1"
$ws.Range("C42").Value = "Right Leg 1"
$ws.Range("D42").Value = "leg_1_r"
$ws.Range("E42").Value = "RightLeg1"
$ws.Range("F42").Value = "R_Leg_1"

$ws.Range("B43").Value = "Right Leg 2"
$ws.Range("C43").Value = "Right Leg 2"
$ws.Range("D43").Value = "leg_2_r"
$ws.Range("E43").Value = "RightLeg1"
$ws.Range("F43").Value = "R_Leg_2"

$ws.Range("B44").Value = "Right Leg 3"
$ws.Range("C44").Value = "Right Leg 3"
$ws.Range("D44").Value = "leg_3_r"
$ws.Range("E44").Value = "RightLeg2"
$ws.Range("F44").Value = "R_Leg_3"

$ws.Range("B45").Value = "Left Leg 1"
$ws.Range("C45").Value = "Left Leg 1"
$ws.Range("D45").Value = "leg_1_l"
$ws.Range("E45").Value = "LeftLeg1"
$ws.Range("F45").Value = "L_Leg_1"

$ws.Range("B46").Value = "Left Leg 2"
$ws.Range("C46").Value = "Left Leg 2"
$ws.Range("D46").Value = "leg_2_l"
$ws.Range("E46").Value = "LeftLeg1"
$ws.Range("F46").Value = "L_Leg_2"

$ws.Range("B47").Value = "Left Leg 3"
$ws.Range("C47").Value = "Left Leg 3"
$ws.Range("D47").Value = "leg_3_l"
$ws.Range("E47").Value = "LeftLeg2"
$ws.Range("F47").Value = "L_Leg_3"

# --- Grow Table18 (the xml-mapped item table) to cover the new rows. ---
$lo2 = $ws.ListObjects.Item(2)
$lo2.Resize($ws.Range("B5:F47"))

# --- Update the view: scroll back to the top and select L17 (matches the
#     saved window/selection state in the edited workbook). ---
$win = $wb.Windows.Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("L17").Select()
